# Logged Week 15 and simulated Week 16
#
# Appends the new weeks' per-play yardage logs (YDS sheet) and per-week
# special-teams logs (ST sheet) to their existing space-separated strings,
# then updates the derived season-to-date totals on OFF, DEF, ST, TURNS and
# PEN with the new running totals that include weeks 15 & 16.

$wb = $excel.ActiveWorkbook

function Append-Tokens($ws, $addr, $tokens) {
    $old = $ws.Range($addr).Value()
    $ws.Range($addr).Value = $old + " " + $tokens
}

# ---------------------------------------------------------------------
# YDS sheet: offensive/defensive rush (R) and pass (P) play-by-play yards
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")
Append-Tokens $ydsWs "B2" "4 5 4 4 4 6 5 8 1 4 19 0 9 4 4 10 3 1 8 8 38 1 9 12 7 -3 4 3 15 4 2 3 4 4 0 8 4 -3 1"
Append-Tokens $ydsWs "B3" "5 -1 19 28 16 45 6 5 8 17 21 20 5 9 7 9 7 34 19 17"
Append-Tokens $ydsWs "C2" "0 2 6 5 2 4 1 3 3 5 4 -1 -1 8 -6 2 30 0 1 0 13"
Append-Tokens $ydsWs "C3" "19 7 6 46 5 1 7 5 9 7 5 12 29 5 6 1 13 10 3 -2"

# ---------------------------------------------------------------------
# ST sheet: kickoff (KO) / punt (PT) depth, return-attempt and return-max
# per-week logs
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")
Append-Tokens $stWs "B4" "64 62"
Append-Tokens $stWs "B5" "47 28"
Append-Tokens $stWs "B6" "31 23 17"
Append-Tokens $stWs "D3" "41 34"
Append-Tokens $stWs "D4" "0 0"
Append-Tokens $stWs "D5" "0 11 0 1 10"

# ---------------------------------------------------------------------
# OFF sheet: updated season totals (Home row = 2, Road row = 3)
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("B2").Value = 4
$offWs.Range("C2").Value = 188
$offWs.Range("D2").Value = 13
$offWs.Range("E2").Value = 10
$offWs.Range("F2").Value = 65
$offWs.Range("G2").Value = 62
$offWs.Range("I2").Value = 13
$offWs.Range("J2").Value = 38
$offWs.Range("L2").Value = 165
$offWs.Range("M2").Value = 100
$offWs.Range("O2").Value = 17
$offWs.Range("P2").Value = 6
$offWs.Range("Q2").Value = 377

$offWs.Range("C3").Value = 122
$offWs.Range("D3").Value = 5
$offWs.Range("E3").Value = 24
$offWs.Range("F3").Value = 77
$offWs.Range("G3").Value = 20
$offWs.Range("H3").Value = 17
$offWs.Range("I3").Value = 33
$offWs.Range("J3").Value = 49
$offWs.Range("N3").Value = 14

# ---------------------------------------------------------------------
# DEF sheet: updated season totals (Home row = 2, Road row = 3)
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("B2").Value = 4
$defWs.Range("C2").Value = 170
$defWs.Range("D2").Value = 12
$defWs.Range("F2").Value = 47
$defWs.Range("G2").Value = 57
$defWs.Range("J2").Value = 22
$defWs.Range("L2").Value = 210
$defWs.Range("M2").Value = 154
$defWs.Range("O2").Value = 23
$defWs.Range("Q2").Value = 378

$defWs.Range("C3").Value = 157
$defWs.Range("E3").Value = 27
$defWs.Range("F3").Value = 100
$defWs.Range("G3").Value = 30
$defWs.Range("H3").Value = 23
$defWs.Range("I3").Value = 53
$defWs.Range("J3").Value = 43
$defWs.Range("N3").Value = 4

# ---------------------------------------------------------------------
# ST sheet: updated season totals (row 2 = "#", row 3 = "TB")
# ---------------------------------------------------------------------
$stWs.Range("B2").Value = 75
$stWs.Range("D2").Value = 45
$stWs.Range("F2").Value = 176
$stWs.Range("G2").Value = 167
$stWs.Range("J2").Value = 68
$stWs.Range("K2").Value = 62
$stWs.Range("B3").Value = 49

# ---------------------------------------------------------------------
# TURNS sheet: updated season totals (Home row = 2, Road row = 3)
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("B2").Value = 2
$turnsWs.Range("D2").Value = 6
$turnsWs.Range("E2").Value = 5
$turnsWs.Range("D3").Value = 4

# ---------------------------------------------------------------------
# PEN sheet: updated season totals
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 13
$penWs.Range("B3").Value = 13
$penWs.Range("D4").Value = 8
